# Applies:
#  1) rows 1480-1510: column R ("backup") goes from blank/inline-empty-string
#     to an explicit numeric 0.
#  2) appends new historical rows 1511-1532 (columns A..Q populated, R left
#     blank, matching the pattern of every prior row in the sheet).
#  3) sheet dimension grows from A1:R1510 to A1:R1532 (handled automatically
#     by Excel once the new cells are written).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastCol = 18  # column R ("backup")

# --- 1) Backfill column R with 0 for existing rows 1480..1510 ---
for ($r = 1480; $r -le 1510; $r++) {
    $ws.Cells.Item($r, $lastCol).Value = 0
}

# --- 2) Append the new rows 1511..1532 ---
# Each inner array is: row, A(Datetime-serial), B(Open), C(High), D(Low),
# E(Close), F(AdjClose), G(Volume), H(Year), I(Month), J(Day), K(Hour),
# L(Minute), M(Second), N(Week), O(isPivot), P(two_line_structure),
# Q(detect_structure). Column R ("backup") is intentionally left blank for
# these new rows, same as every row before the 1480-1510 backfill above.
$newRows = @(
    @(1511, 45555, 4392, 4414.9501953125, 4360, 4378.2998046875, 4378.2998046875, 737967, 2024, 9, 20, 0, 0, 0, 38, 0, 0, 2),
    @(1512, 45558, 4385.5, 4433, 4330, 4412.7998046875, 4412.7998046875, 415628, 2024, 9, 23, 0, 0, 0, 39, 0, 0, 0),
    @(1513, 45559, 4422, 4534.2998046875, 4390, 4479.14990234375, 4479.14990234375, 558134, 2024, 9, 24, 0, 0, 0, 39, 0, 0, 0),
    @(1514, 45560, 4489.89990234375, 4505.7001953125, 4425.0498046875, 4471.2001953125, 4471.2001953125, 261274, 2024, 9, 25, 0, 0, 0, 39, 0, 0, 0),
    @(1515, 45561, 4490.39990234375, 4541, 4416, 4527.5498046875, 4527.5498046875, 268104, 2024, 9, 26, 0, 0, 0, 39, 1, 0, 0),
    @(1516, 45562, 4525, 4525, 4411, 4420.5498046875, 4420.5498046875, 496527, 2024, 9, 27, 0, 0, 0, 39, 0, 0, 0),
    @(1517, 45565, 4420, 4420, 4275.5498046875, 4299.2998046875, 4299.2998046875, 398612, 2024, 9, 30, 0, 0, 0, 40, 0, 0, 0),
    @(1518, 45566, 4318, 4390.7001953125, 4301.75, 4365, 4365, 335988, 2024, 10, 1, 0, 0, 0, 40, 0, 0, 0),
    @(1519, 45568, 4292, 4350, 4201.0498046875, 4225.64990234375, 4225.64990234375, 231145, 2024, 10, 3, 0, 0, 0, 40, 0, 0, 0),
    @(1520, 45569, 4225, 4289.64990234375, 4184.64990234375, 4219.4501953125, 4219.4501953125, 285684, 2024, 10, 4, 0, 0, 0, 40, 0, 0, 0),
    @(1521, 45572, 4203, 4237.9501953125, 4102.25, 4132.7998046875, 4132.7998046875, 483253, 2024, 10, 7, 0, 0, 0, 41, 0, 0, 0),
    @(1522, 45573, 4100, 4202, 4100, 4189.4501953125, 4189.4501953125, 371923, 2024, 10, 8, 0, 0, 0, 41, 2, 0, 0),
    @(1523, 45574, 4249.89990234375, 4388.9501953125, 4239.35009765625, 4369.7998046875, 4369.7998046875, 692046, 2024, 10, 9, 0, 0, 0, 41, 0, 0, 0),
    @(1524, 45575, 4360.5, 4443.5498046875, 4352.25, 4410.4501953125, 4410.4501953125, 385531, 2024, 10, 10, 0, 0, 0, 41, 0, 0, 0),
    @(1525, 45576, 4411, 4455, 4345, 4449.39990234375, 4449.39990234375, 717931, 2024, 10, 11, 0, 0, 0, 41, 0, 0, 0),
    @(1526, 45579, 4476.9501953125, 4502, 4414, 4482.0498046875, 4482.0498046875, 497069, 2024, 10, 14, 0, 0, 0, 42, 0, 2, 2),
    @(1527, 45580, 4499.89990234375, 4604.25, 4465.14990234375, 4555.10009765625, 4555.10009765625, 1566800, 2024, 10, 15, 0, 0, 0, 42, 0, 0, 0),
    @(1528, 45581, 4560, 4864, 4560, 4827.89990234375, 4827.89990234375, 3365030, 2024, 10, 16, 0, 0, 0, 42, 0, 0, 0),
    @(1529, 45582, 4824.0498046875, 4824.0498046875, 4690.2998046875, 4728.64990234375, 4728.64990234375, 670142, 2024, 10, 17, 0, 0, 0, 42, 0, 0, 0),
    @(1530, 45583, 4698.9501953125, 4750.2001953125, 4605, 4709.75, 4709.75, 343633, 2024, 10, 18, 0, 0, 0, 42, 0, 0, 0),
    @(1531, 45586, 4740.89990234375, 4800, 4590.0498046875, 4602.9501953125, 4602.9501953125, 607174, 2024, 10, 21, 0, 0, 0, 43, 0, 0, 0),
    @(1532, 45587, 4602.9501953125, 4647.89990234375, 4490.5498046875, 4510.85009765625, 4510.85009765625, 567025, 2024, 10, 22, 0, 0, 0, 43, 0, 0, 0)
)

foreach ($row in $newRows) {
    $r = $row[0]

    # Column A carries the same "yyyy-mm-dd hh:mm:ss" datetime style (s="2")
    # as every other row in the sheet.
    $cellA = $ws.Cells.Item($r, 1)
    $cellA.NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $cellA.Value = $row[1]

    for ($col = 2; $col -le 17; $col++) {
        $ws.Cells.Item($r, $col).Value = $row[$col]
    }

    # Column R ("backup") stays blank for the newly appended rows.
}

Write-Host "done"
